$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 167, shifting existing rows 167-235 down to 168-236.
$ws.Rows.Item(167).Insert()

# Populate the new row 167 with the new record (duplicating row-166 layout/labels,
# which are constant across this worksheet's rows).
$ws.Range("A167").Value = 7
$ws.Range("B167").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C167").Value = "Ñuble"
$ws.Range("D167").Value = 44755
$ws.Range("D167").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E167").Value = 16
$ws.Range("F167").Value = 100112043
$ws.Range("G167").Value = "Pepino ensalada"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 100
$ws.Range("K167").Value = 19000
$ws.Range("L167").Value = 20000
$ws.Range("M167").Value = 19500
$ws.Range("N167").Value = "$/caja 60 unidades"
$ws.Range("O167").Value = "Región de Arica y Parinacota"
$ws.Range("P167").Value = 325
$ws.Range("Q167").Value = 60
$ws.Range("R167").Value = "Hortaliza"
